$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices formatted as plain text (e.g. "68.005.06"); force text
# storage before assigning so Excel does not reinterpret them as numbers.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

# Apply the refreshed price / 1h-volume figures
$ws.Range('D2').Value = '68.005.06'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '3.524.42'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '602.43'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').Value = '182.02'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.525.75'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('D9').Value = '0.597'
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('E10').Value = '  +6.18%  '
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('D12').Value = '0.440'
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('D13').Value = '4.134.11'
$ws.Range('D14').Value = '32.40'
$ws.Range('E14').Value = '  +10.54%  '
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').Value = '67.956.32'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').Value = '3.526.26'
$ws.Range('E18').Value = '  -1.43%  '
$ws.Range('D19').Value = '6.38'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('D20').Value = '14.52'
$ws.Range('E20').Value = '  +2.49%  '
$ws.Range('D21').Value = '401.55'
$ws.Range('E21').Value = '  +1.33%  '
$ws.Range('D22').Value = '8.02'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').Value = '74.09'
$ws.Range('E23').Value = '  +1.34%  '
$ws.Range('D24').Value = '0.546'
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('D25').Value = '0.997'
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('D28').Value = '10.58'
$ws.Range('E28').Value = '  +2.92%  '
$ws.Range('E29').Value = '  -2.39%  '
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('D32').Value = '1.45'
$ws.Range('E32').Value = '  -1.05%  '
$ws.Range('D33').Value = '2.09'
$ws.Range('E33').Value = '  +1.05%  '
$ws.Range('D34').Value = '24.01'
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('D35').Value = '7.54'
$ws.Range('E35').Value = '  +1.83%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').Value = '  -2.19%  '
$ws.Range('D38').Value = '163.19'
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').Value = '0.884'
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('D40').Value = '1.93'
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('E41').Value = '  +8.45%  '
$ws.Range('D42').Value = '7.00'
$ws.Range('E42').Value = '  -0.73%  '
$ws.Range('D43').Value = '4.73'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('D44').Value = '2.888.50'
$ws.Range('E44').Value = '  +1.68%  '
$ws.Range('D45').Value = '26.55'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').Value = '0.0737'
$ws.Range('E46').Value = '  -2.07%  '
$ws.Range('D47').Value = '26.90'
$ws.Range('E47').Value = '  -1.27%  '
$ws.Range('D48').Value = '42.57'
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('D49').Value = '351.84'
$ws.Range('E49').Value = '  +3.81%  '
$ws.Range('D50').Value = '0.0306'
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('E51').Value = '  -1.06%  '
